{"js": "// Apply the commit's edits to the DT_Floriane_RENAUD document.\n//\n// 1) Title (first table, single cell) -> uppercase.\n// 2) \"Domaine\" row (table 2, row 3) -> add a space, insert a new\n//    \"Suite Office\" line, and change \"5S\" to \"5$S\".\n// 3) \"Normes\" row (table 2, row 4) -> consolidate the CENELEC / ISO\n//    lines.\n// 4) \"Logiciels\" row (table 2, row 5) -> was empty, now has two\n//    lines of text.\n// 5) Final job-title paragraph -> new text.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\n// --- 1) Title table ------------------------------------------------\nconst titleTable = tables.items[0];\nconst titleCell = titleTable.getCell(0, 0);\ntitleCell.body.load(\"text\");\nawait context.sync();\n\nconst oldTitle = \"Manager en Maitrise des Risques Industriels\";\nconst newTitle = \"MANAGER EN MAITRISE DES RISQUES INDUSTRIELS\";\nif (titleCell.body.text.trim() === oldTitle) {\n  // Keep the existing run formatting (sz 28, centered) by replacing the\n  // text in place rather than clearing the cell.\n  const hits = titleCell.body.search(oldTitle, { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  hits.items[0].insertText(newTitle, \"Replace\");\n  await context.sync();\n}\n\n// --- 2) \"Domaine de comp\u00e9tences\" skills table -----------------------\nconst skillsTable = tables.items[1];\n\n// Row index 2 -> \"Domaine\" row, second column.\nconst domaineCell = skillsTable.getCell(2, 1);\ndomaineCell.body.clear();\ndomaineCell.body.insertText(\n  \"Normes (CENELEC 50128/ 50126/ 50129, ISO 9001 & 140001, ISO/IEC 12207, MIL-STD-498)\\u000b\" +\n  \"Suite Office (Word, Excel, Powerpoint)\\u000b\" +\n  \"Ingenierie pedagogique\\u000b\" +\n  \"Gestion documentaire\\u000b\" +\n  \"Methodes 5$S, 8D, Ishikawa\",\n  \"Replace\"\n);\n\n// Row index 3 -> \"Normes\" row, second column.\nconst normesCell = skillsTable.getCell(3, 1);\nnormesCell.body.clear();\nnormesCell.body.insertText(\n  \"CENELEC 50128/50126/50129\\u000b\" +\n  \"ISO 9001 & 140001\\u000b\" +\n  \"ISO/IEC 12207\\u000b\" +\n  \"MIL-STD-498\",\n  \"Replace\"\n);\n\n// Row index 4 -> \"Logiciels\" row, second column (previously empty).\nconst logicielsCell = skillsTable.getCell(4, 1);\nlogicielsCell.body.clear();\nlogicielsCell.body.insertText(\n  \"Suite Office (Word, Excel, Powerpoint)\\u000b\" +\n  \"Gestion documentaire\",\n  \"Replace\"\n);\n\nawait context.sync();\n\n// --- 3) Final job-title paragraph -----------------------------------\nconst body = context.document.body;\nconst oldJobTitle = \"Ingenieur en Surete de fonctionnement\";\nconst newJobTitle = \"Ingenieur d'etudes & qualite\";\nconst jobHits = body.search(oldJobTitle, { matchCase: true });\njobHits.load(\"items\");\nawait context.sync();\n\nif (jobHits.items.length > 0) {\n  jobHits.items[0].insertText(newJobTitle, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Apply the commit's edits to the DT_Floriane_RENAUD document.\n#\n# 1) Title (first table, single cell) -> uppercase.\n# 2) \"Domaine\" row (table 2, row 3) -> add a space, insert a new\n#    \"Suite Office\" line, and change \"5S\" to \"5$S\".\n# 3) \"Normes\" row (table 2, row 4) -> consolidate the CENELEC / ISO\n#    lines.\n# 4) \"Logiciels\" row (table 2, row 5) -> was empty, now has two\n#    lines of text.\n# 5) Final job-title paragraph -> new text.\n\n$d = $word.ActiveDocument\n$lineBreak = [char]11\n\n# --- 1) Title table --------------------------------------------------\n$titleTable = $d.Tables.Item(1)\n$titleCell = $titleTable.Cell(1, 1)\n$titleRange = $titleCell.Range\n$titleRange.Find.Execute(\"Manager en Maitrise des Risques Industriels\")\nif ($titleRange.Find.Found) {\n  # Assign straight to Range.Text (not Find's ReplaceWith) so the run's\n  # existing formatting (size 28, centered) is kept and no smart-quote /\n  # autocorrect substitution is applied.\n  $titleRange.Text = \"MANAGER EN MAITRISE DES RISQUES INDUSTRIELS\"\n}\n\n# --- 2) \"Domaine de comp\u00e9tences\" skills table -------------------------\n$skillsTable = $d.Tables.Item(2)\n\n# Row 3 -> \"Domaine\" row, second column.\n$domaineCell = $skillsTable.Cell(3, 2)\n$domaineCell.Range.Text = (\n  \"Normes (CENELEC 50128/ 50126/ 50129, ISO 9001 & 140001, ISO/IEC 12207, MIL-STD-498)\" + $lineBreak +\n  \"Suite Office (Word, Excel, Powerpoint)\" + $lineBreak +\n  \"Ingenierie pedagogique\" + $lineBreak +\n  \"Gestion documentaire\" + $lineBreak +\n  \"Methodes 5`$S, 8D, Ishikawa\"\n)\n\n# Row 4 -> \"Normes\" row, second column.\n$normesCell = $skillsTable.Cell(4, 2)\n$normesCell.Range.Text = (\n  \"CENELEC 50128/50126/50129\" + $lineBreak +\n  \"ISO 9001 & 140001\" + $lineBreak +\n  \"ISO/IEC 12207\" + $lineBreak +\n  \"MIL-STD-498\"\n)\n\n# Row 5 -> \"Logiciels\" row, second column (previously empty).\n$logicielsCell = $skillsTable.Cell(5, 2)\n$logicielsCell.Range.Text = (\n  \"Suite Office (Word, Excel, Powerpoint)\" + $lineBreak +\n  \"Gestion documentaire\"\n)\n\n# --- 3) Final job-title paragraph -------------------------------------\n$jobRange = $d.Content\n$jobRange.Find.Execute(\"Ingenieur en Surete de fonctionnement\")\nif ($jobRange.Find.Found) {\n  $jobRange.Text = \"Ingenieur d'etudes & qualite\"\n}\n"}
